# Applies the 2024-05-26 cryptos-list refresh (prices + 1h volume %, plus two
# rows whose ranking swapped: NEARProtocol/PancakeSwap and Stacks/Bittensor).
#
# The Price/Volume columns are plain text cells (e.g. "68.923.45", "  +0.10%  "),
# not numbers. Assigning a numeric-looking string straight to Range.Value lets
# Excel auto-coerce it to a real number (losing the original text formatting),
# so values that would otherwise parse as a float are written with a leading
# apostrophe ('=text-literal prefix) to force them to stay text, matching the
# source data's own representation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range('D2').Value = '68.923.45'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '3.738.25'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''601.73'
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').Value = '''165.78'
$ws.Range('E6').Value = '  -1.89%  '
$ws.Range('D7').Value = '3.733.31'
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +0.92%  '
$ws.Range('E10').Value = '  +3.70%  '
$ws.Range('D11').Value = '''6.38'
$ws.Range('E11').Value = '  +0.80%  '
$ws.Range('D12').Value = '''0.459'
$ws.Range('E12').Value = '  -0.34%  '
$ws.Range('D13').Value = '''37.76'
$ws.Range('E13').Value = '  -1.15%  '
$ws.Range('D14').Value = '''0.0000247'
$ws.Range('E14').Value = '  +0.54%  '
$ws.Range('D15').Value = '4.365.67'
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('D16').Value = '3.740.01'
$ws.Range('E16').Value = '  +0.25%  '
$ws.Range('D17').Value = '69.008.02'
$ws.Range('E17').Value = '  +0.18%  '
$ws.Range('D18').Value = '''7.43'
$ws.Range('E18').Value = '  +2.05%  '
$ws.Range('E19').Value = '  +3.27%  '
$ws.Range('E20').Value = '  -2.12%  '
$ws.Range('D21').Value = '''11.19'
$ws.Range('E21').Value = '  +5.13%  '
$ws.Range('D22').Value = '''492.33'
$ws.Range('D23').Value = '''0.725'
$ws.Range('E23').Value = '  -0.32%  '
$ws.Range('D24').Value = '''84.66'
$ws.Range('E24').Value = '  -0.63%  '
$ws.Range('E25').Value = '  +2.37%  '
$ws.Range('E26').Value = '  -1.86%  '
$ws.Range('E27').Value = '  -0.70%  '
$ws.Range('D28').Value = '''10.07'
$ws.Range('E28').Value = '  -0.91%  '
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '''2.97'
$ws.Range('E30').Value = '  -0.17%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').Value = '''8.19'
$ws.Range('E31').Value = '  +3.01%  '
$ws.Range('E32').Value = '  -5.01%  '
$ws.Range('D33').Value = '''31.58'
$ws.Range('E33').Value = '  -0.24%  '
$ws.Range('D34').Value = '3.884.91'
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').Value = '3.675.11'
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('E36').Value = '  -0.17%  '
$ws.Range('E37').Value = '  +1.51%  '
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('E39').Value = '  +5.21%  '
$ws.Range('D40').Value = '''3.17'
$ws.Range('E40').Value = '  +9.90%  '
$ws.Range('D41').Value = '''1.00'
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('E42').Value = '  -0.11%  '
$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D43').Value = '''427.02'
$ws.Range('E43').Value = '  -3.11%  '
$ws.Range('D44').Value = '''48.59'
$ws.Range('E44').Value = '  -0.78%  '
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').Value = '''1.99'
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('D46').Value = '''8.44'
$ws.Range('E46').Value = '  -0.62%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').Value = '''39.95'
$ws.Range('E48').Value = '  -1.85%  '
$ws.Range('D49').Value = '''141.02'
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('D50').Value = '''1.29'
$ws.Range('E50').Value = '  +6.71%  '
$ws.Range('D51').Value = '2.777.18'
$ws.Range('E51').Value = '  -0.27%  '
